$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.686300000000004
$ws.Range("C8").Value = -14.4875
$ws.Range("C10").Value = -13.25029999999999
$ws.Range("D11").Value = -7.012699999999994
$ws.Range("C12").Value = -13.2601
$ws.Range("D12").Value = -8.588400000000002
$ws.Range("D15").Value = -8.546199999999997
$ws.Range("D17").Value = -8.308199999999992
$ws.Range("C18").Value = -10.8224
$ws.Range("C25").Value = -10.8538
$ws.Range("D26").Value = -7.251700000000003
$ws.Range("D27").Value = -8.439400000000001
$ws.Range("D28").Value = -8.553799999999999
$ws.Range("D32").Value = -6.309499999999996
$ws.Range("C37").Value = -13.59340000000001
$ws.Range("D37").Value = -7.077000000000002
$ws.Range("D41").Value = -8.311699999999991
$ws.Range("D47").Value = -7.848000000000002
$ws.Range("D51").Value = -8.6119
$ws.Range("C55").Value = -12.9402
$ws.Range("D65").Value = -7.924800000000005
$ws.Range("C68").Value = -10.781
$ws.Range("D73").Value = -7.995599999999999
$ws.Range("C77").Value = -13.06610000000001
$ws.Range("C78").Value = -12.65800000000001
$ws.Range("C79").Value = -12.30890000000001
$ws.Range("C80").Value = -12.72160000000001
$ws.Range("C81").Value = -14.6933
$ws.Range("C82").Value = -11.07769999999999
$ws.Range("C84").Value = -13.0872
$ws.Range("D84").Value = -8.515400000000001
$ws.Range("D85").Value = -8.223500000000001
$ws.Range("D89").Value = -8.402599999999993
$ws.Range("D93").Value = -6.450499999999992
$ws.Range("D95").Value = -7.626900000000001
$ws.Range("D98").Value = -7.0448
$ws.Range("D99").Value = -8.291100000000005
$ws.Range("C101").Value = -13.22459999999999
$ws.Range("D101").Value = -7.7341
$ws.Range("C102").Value = -12.62930000000001
$ws.Range("D102").Value = -7.072899999999996
